# "Europa" window edit: replace the Belgium/Europe stats block with
# an "Inform about Asia" summary block, and drop the trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: country/title -> section header, and its value -> big number (as text)
$ws.Range("A1").Value = "Інформація по Азії"
$ws.Range("B1").Value = "'81,532,629"
$ws.Range("B1").ClearFormats()

# Row 2: last-update label -> "total cases" label, value -> total cases (as text)
$ws.Range("A2").Value = "Всього випадків"
$ws.Range("B2").Value = "'1,204,513"
$ws.Range("B2").ClearFormats()

# Row 3: total-cases label -> total-deaths label, numeric value -> text value
$ws.Range("A3").Value = "Загальна кількість смертей"
$ws.Range("B3").Value = "'78,786,289"
$ws.Range("B3").ClearFormats()

# Row 4: keep only the "Число одуживших" label (from the old row 6); no value.
$ws.Range("A4").Value = "Число одуживших"
$ws.Range("B4").ClearContents()

# Remove what used to be rows 5, 6 and 7 entirely.
$ws.Range("A5:B7").Clear()
